$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Cells.Item(28, 8).Value = 1487.1111  # H28: was 776.125
$ws.Cells.Item(28, 9).Value = 1045.1875  # I28: was 824.86664
$ws.Cells.Item(28, 10).Value = 5022.5  # J28: was 45
$ws.Cells.Item(28, 11).Value = 1045.1875  # K28: was 824.86664
$ws.Cells.Item(28, 12).Value = 5022.5  # L28: was 45
$ws.Cells.Item(28, 13).Value = -560.1875  # M28: was -339.86664
$ws.Cells.Item(28, 14).Value = -5992.5  # N28: was -1015
$ws.Cells.Item(32, 8).Value = 7816.4546  # H32: was 7870.091
$ws.Cells.Item(32, 9).Value = 6917.909  # I32: was 7025.1816
$ws.Cells.Item(32, 11).Value = 6917.909  # K32: was 7025.1816
$ws.Cells.Item(32, 13).Value = -6591.909  # M32: was -6699.1816
$ws.Cells.Item(80, 8).Value = 10970.685  # H80: was 9928.096
$ws.Cells.Item(80, 9).Value = 8597.916999999999  # I80: was 7937.3076
$ws.Cells.Item(80, 10).Value = 15038.286  # J80: was 13163.125
$ws.Cells.Item(80, 11).Value = 25793.751  # K80: was 23811.9228
$ws.Cells.Item(80, 12).Value = 45114.858  # L80: was 39489.375
$ws.Cells.Item(80, 13).Value = -24795.751  # M80: was -22813.9228
$ws.Cells.Item(80, 14).Value = -47110.858  # N80: was -41485.375
$ws.Cells.Item(83, 8).Value = 10970.685  # H83: was 9928.096
$ws.Cells.Item(83, 9).Value = 8597.916999999999  # I83: was 7937.3076
$ws.Cells.Item(83, 10).Value = 15038.286  # J83: was 13163.125
$ws.Cells.Item(83, 11).Value = 77381.253  # K83: was 71435.7684
$ws.Cells.Item(83, 12).Value = 135344.574  # L83: was 118468.125
$ws.Cells.Item(83, 13).Value = -72389.253  # M83: was -66443.7684
$ws.Cells.Item(83, 14).Value = -145328.574  # N83: was -128452.125
$ws.Cells.Item(94, 8).Value = 998.3333  # H94: was 1000
$ws.Cells.Item(94, 9).Value = 998.3333  # I94: was 1000
$ws.Cells.Item(94, 11).Value = 998.3333  # K94: was 1000
$ws.Cells.Item(94, 13).Value = -547.3333  # M94: was -549
$ws.Cells.Item(118, 8).Value = 496.33334  # H118: was 494.5
$ws.Cells.Item(118, 9).Value = 499.5  # I118: was 496
$ws.Cells.Item(118, 11).Value = 1498.5  # K118: was 1488
$ws.Cells.Item(118, 13).Value = 158.5  # M118: was 169
$ws.Cells.Item(137, 8).Value = 8639.706  # H137: was 8995.4375
$ws.Cells.Item(137, 9).Value = 9215.23  # I137: was 9737.5
$ws.Cells.Item(137, 11).Value = 27645.69  # K137: was 29212.5
$ws.Cells.Item(137, 13).Value = -25095.69  # M137: was -26662.5
$ws.Cells.Item(138, 8).Value = 5076.8096  # H138: was 4851.617
$ws.Cells.Item(138, 10).Value = 4663.1  # J138: was 4419.8
$ws.Cells.Item(138, 12).Value = 13989.3  # L138: was 13259.4
$ws.Cells.Item(138, 14).Value = -24269.3  # N138: was -23539.4
$ws.Cells.Item(141, 8).Value = 2169.8333  # H141: was 2045.5714
$ws.Cells.Item(141, 9).Value = 2003.3529  # I141: was 1897.85
$ws.Cells.Item(141, 11).Value = 6010.0587  # K141: was 5693.549999999999
$ws.Cells.Item(141, 13).Value = -830.0587000000005  # M141: was -513.5499999999993

# --- Sheet: ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Cells.Item(32, 8).Value = 13641.612  # H32: was 12685.438
$ws.Cells.Item(32, 9).Value = 11770.328  # I32: was 10896.1045
$ws.Cells.Item(32, 11).Value = 11770.328  # K32: was 10896.1045
$ws.Cells.Item(32, 13).Value = -11483.328  # M32: was -10609.1045
$ws.Cells.Item(45, 8).Value = 103314.05  # H45: was 108661.52
$ws.Cells.Item(45, 9).Value = 169365.5  # I45: was 184606.73
$ws.Cells.Item(45, 11).Value = 169365.5  # K45: was 184606.73
$ws.Cells.Item(45, 13).Value = -168988.5  # M45: was -184229.73
$ws.Cells.Item(61, 8).Value = 2455.3333  # H61: was 2455.7222
$ws.Cells.Item(61, 9).Value = 2455.3333  # I61: was 2553.2354
$ws.Cells.Item(61, 10).Value = 0  # J61: was 798
$ws.Cells.Item(61, 11).Value = 2455.3333  # K61: was 2553.2354
$ws.Cells.Item(61, 12).Value = 0  # L61: was 798
$ws.Cells.Item(61, 13).Value = -2243.3333  # M61: was -2341.2354
$ws.Cells.Item(61, 14).ClearContents()  # N61: was -1222
$ws.Cells.Item(63, 8).Value = 2385.8462  # H63: was 2337
$ws.Cells.Item(63, 9).Value = 2573.6365  # I63: was 2412.4614
$ws.Cells.Item(63, 10).Value = 1353  # J63: was 1356
$ws.Cells.Item(63, 11).Value = 2573.6365  # K63: was 2412.4614
$ws.Cells.Item(63, 12).Value = 1353  # L63: was 1356
$ws.Cells.Item(63, 13).Value = -1887.6365  # M63: was -1726.4614
$ws.Cells.Item(63, 14).Value = -2725  # N63: was -2728
$ws.Cells.Item(66, 8).Value = 2385.8462  # H66: was 2337
$ws.Cells.Item(66, 9).Value = 2573.6365  # I66: was 2412.4614
$ws.Cells.Item(66, 10).Value = 1353  # J66: was 1356
$ws.Cells.Item(66, 11).Value = 12868.1825  # K66: was 12062.307
$ws.Cells.Item(66, 12).Value = 6765  # L66: was 6780
$ws.Cells.Item(66, 13).Value = -9436.182500000001  # M66: was -8630.307000000001
$ws.Cells.Item(66, 14).Value = -13629  # N66: was -13644
$ws.Cells.Item(74, 8).Value = 889.9231  # H74: was 944.9091
$ws.Cells.Item(74, 9).Value = 839.0833  # I74: was 889.4
$ws.Cells.Item(74, 11).Value = 839.0833  # K74: was 889.4
$ws.Cells.Item(74, 13).Value = 34.91669999999999  # M74: was -15.39999999999998
$ws.Cells.Item(77, 8).Value = 889.9231  # H77: was 944.9091
$ws.Cells.Item(77, 9).Value = 839.0833  # I77: was 889.4
$ws.Cells.Item(77, 11).Value = 4195.4165  # K77: was 4447
$ws.Cells.Item(77, 13).Value = 172.5834999999997  # M77: was -79
$ws.Cells.Item(132, 8).Value = 2943.5  # H132: was 2994.3572
$ws.Cells.Item(132, 9).Value = 2943.5  # I132: was 2994.3572
$ws.Cells.Item(132, 11).Value = 8830.5  # K132: was 8983.071599999999
$ws.Cells.Item(132, 13).Value = -6300.5  # M132: was -6453.071599999999
$ws.Cells.Item(136, 8).Value = 2455.3333  # H136: was 2455.7222
$ws.Cells.Item(136, 9).Value = 2455.3333  # I136: was 2553.2354
$ws.Cells.Item(136, 10).Value = 0  # J136: was 798
$ws.Cells.Item(136, 11).Value = 7365.999899999999  # K136: was 7659.706200000001
$ws.Cells.Item(136, 12).Value = 0  # L136: was 2394
$ws.Cells.Item(136, 13).Value = -4815.999899999999  # M136: was -5109.706200000001
$ws.Cells.Item(136, 14).ClearContents()  # N136: was -7494
$ws.Cells.Item(139, 8).Value = 71498.5  # H139: was 55749.5
$ws.Cells.Item(139, 10).Value = 71498.5  # J139: was 55749.5
$ws.Cells.Item(139, 12).Value = 71498.5  # L139: was 55749.5
$ws.Cells.Item(139, 14).Value = -81778.5  # N139: was -66029.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Cells.Item(99, 8).Value = 4999  # H99: was 1561.2
$ws.Cells.Item(99, 9).Value = 0  # I99: was 701.75
$ws.Cells.Item(99, 11).Value = 0  # K99: was 701.75
$ws.Cells.Item(99, 13).ClearContents()  # M99: was 796.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Cells.Item(16, 8).Value = 3749.9  # H16: was 3937.375
$ws.Cells.Item(16, 9).Value = 3428.5715  # I16: was 3600
$ws.Cells.Item(16, 11).Value = 3428.5715  # K16: was 3600
$ws.Cells.Item(16, 13).Value = -3141.5715  # M16: was -3313
$ws.Cells.Item(31, 8).Value = 5001.2666  # H31: was 5316.077
$ws.Cells.Item(31, 9).Value = 3281.4285  # I31: was 3076.125
$ws.Cells.Item(31, 10).Value = 6506.125  # J31: was 8900
$ws.Cells.Item(31, 11).Value = 3281.4285  # K31: was 3076.125
$ws.Cells.Item(31, 12).Value = 6506.125  # L31: was 8900
$ws.Cells.Item(31, 13).Value = -2986.4285  # M31: was -2781.125
$ws.Cells.Item(31, 14).Value = -7096.125  # N31: was -9490
$ws.Cells.Item(34, 8).Value = 5001.2666  # H34: was 5316.077
$ws.Cells.Item(34, 9).Value = 3281.4285  # I34: was 3076.125
$ws.Cells.Item(34, 10).Value = 6506.125  # J34: was 8900
$ws.Cells.Item(34, 11).Value = 3281.4285  # K34: was 3076.125
$ws.Cells.Item(34, 12).Value = 6506.125  # L34: was 8900
$ws.Cells.Item(34, 13).Value = -3079.4285  # M34: was -2874.125
$ws.Cells.Item(34, 14).Value = -6910.125  # N34: was -9304
$ws.Cells.Item(113, 8).Value = 3749.9  # H113: was 3937.375
$ws.Cells.Item(113, 9).Value = 3428.5715  # I113: was 3600
$ws.Cells.Item(113, 11).Value = 3428.5715  # K113: was 3600
$ws.Cells.Item(113, 13).Value = -1258.5715  # M113: was -1430
$ws.Cells.Item(134, 8).Value = 2421.9216  # H134: was 2301.8909
$ws.Cells.Item(134, 9).Value = 1325.6744  # I134: was 1300.9556
$ws.Cells.Item(134, 10).Value = 8314.25  # J134: was 6806.1
$ws.Cells.Item(134, 11).Value = 3977.023200000001  # K134: was 3902.8668
$ws.Cells.Item(134, 12).Value = 24942.75  # L134: was 20418.3
$ws.Cells.Item(134, 13).Value = -1442.023200000001  # M134: was -1367.8668
$ws.Cells.Item(134, 14).Value = -30012.75  # N134: was -25488.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Cells.Item(4, 8).Value = 19026050  # H4: was 17728854
$ws.Cells.Item(4, 9).Value = 15943168  # I4: was 14264990
$ws.Cells.Item(4, 10).Value = 34000050  # J4: was 39666668
$ws.Cells.Item(4, 11).Value = 47829504  # K4: was 42794970
$ws.Cells.Item(4, 12).Value = 102000150  # L4: was 119000004
$ws.Cells.Item(4, 13).Value = -47829392  # M4: was -42794858
$ws.Cells.Item(4, 14).Value = -102000374  # N4: was -119000228
$ws.Cells.Item(7, 8).Value = 566.3333  # H7: was 524.75
$ws.Cells.Item(7, 10).Value = 549.5  # J7: was 499.66666
$ws.Cells.Item(7, 12).Value = 1648.5  # L7: was 1498.99998
$ws.Cells.Item(7, 14).Value = -1872.5  # N7: was -1722.99998
$ws.Cells.Item(12, 8).Value = 995.1111  # H12: was 1119.25
$ws.Cells.Item(12, 10).Value = 995.1111  # J12: was 1119.25
$ws.Cells.Item(12, 12).Value = 2985.3333  # L12: was 3357.75
$ws.Cells.Item(12, 14).Value = -3331.3333  # N12: was -3703.75
$ws.Cells.Item(34, 8).Value = 722.8889  # H34: was 757.1429000000001
$ws.Cells.Item(34, 10).Value = 941.2  # J34: was 1166.6666
$ws.Cells.Item(34, 12).Value = 2823.6  # L34: was 3499.9998
$ws.Cells.Item(34, 14).Value = -2991.6  # N34: was -3667.9998
$ws.Cells.Item(92, 8).Value = 254.57143  # H92: was 278.625
$ws.Cells.Item(92, 10).Value = 270  # J92: was 299.5
$ws.Cells.Item(92, 12).Value = 810  # L92: was 898.5
$ws.Cells.Item(92, 14).Value = -3306  # N92: was -3394.5
$ws.Cells.Item(132, 8).Value = 2071.3572  # H132: was 2499.9
$ws.Cells.Item(132, 10).Value = 2153.8462  # J132: was 2666.6667
$ws.Cells.Item(132, 12).Value = 19384.6158  # L132: was 24000.0003
$ws.Cells.Item(132, 14).Value = -24444.6158  # N132: was -29060.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Cells.Item(11, 8).Value = 126000984  # H11: was 118589220
$ws.Cells.Item(11, 9).Value = 182273800  # I11: was 167084380
$ws.Cells.Item(11, 11).Value = 182273800  # K11: was 167084380
$ws.Cells.Item(11, 13).Value = -182273661  # M11: was -167084241
$ws.Cells.Item(12, 8).Value = 0  # H12: was 7999
$ws.Cells.Item(12, 10).Value = 0  # J12: was 7999
$ws.Cells.Item(12, 12).Value = 0  # L12: was 7999
$ws.Cells.Item(12, 14).ClearContents()  # N12: was -8279
$ws.Cells.Item(14, 8).Value = 6483.5  # H14: was 6057.2856
$ws.Cells.Item(14, 9).Value = 4444  # I14: was 4129.3335
$ws.Cells.Item(14, 11).Value = 4444  # K14: was 4129.3335
$ws.Cells.Item(14, 13).Value = -4276  # M14: was -3961.3335

# --- Sheet: LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Cells.Item(16, 8).Value = 611.25  # H16: was 580
$ws.Cells.Item(16, 10).Value = 795.6667  # J16: was 679.25
$ws.Cells.Item(16, 12).Value = 795.6667  # L16: was 679.25
$ws.Cells.Item(16, 14).Value = -1135.6667  # N16: was -1019.25
$ws.Cells.Item(20, 8).Value = 25090000  # H20: was 33413334
$ws.Cells.Item(20, 10).Value = 25090000  # J20: was 33413334
$ws.Cells.Item(20, 12).Value = 25090000  # L20: was 33413334
$ws.Cells.Item(20, 14).Value = -25090452  # N20: was -33413786
$ws.Cells.Item(117, 8).Value = 29997.5  # H117: was 44000
$ws.Cells.Item(117, 9).Value = 29997.5  # I117: was 30000
$ws.Cells.Item(117, 10).Value = 0  # J117: was 58000
$ws.Cells.Item(117, 11).Value = 29997.5  # K117: was 30000
$ws.Cells.Item(117, 12).Value = 0  # L117: was 58000
$ws.Cells.Item(117, 13).Value = -25408.5  # M117: was -25411
$ws.Cells.Item(117, 14).ClearContents()  # N117: was -67178

# --- Sheet: WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Cells.Item(76, 8).Value = 17999  # H76: was 18000
$ws.Cells.Item(76, 10).Value = 17999  # J76: was 18000
$ws.Cells.Item(76, 12).Value = 17999  # L76: was 18000
$ws.Cells.Item(76, 14).Value = -18629  # N76: was -18630
$ws.Cells.Item(79, 8).Value = 17999  # H79: was 18000
$ws.Cells.Item(79, 10).Value = 17999  # J79: was 18000
$ws.Cells.Item(79, 12).Value = 17999  # L79: was 18000
$ws.Cells.Item(79, 14).Value = -20183  # N79: was -20184
$ws.Cells.Item(113, 8).Value = 2509.842  # H113: was 2755.7058
$ws.Cells.Item(113, 9).Value = 1534.9231  # I113: was 1737.6364
$ws.Cells.Item(113, 11).Value = 4604.7693  # K113: was 5212.9092
$ws.Cells.Item(113, 13).Value = -2434.7693  # M113: was -3042.9092
$ws.Cells.Item(118, 8).Value = 40000  # H118: was 0
$ws.Cells.Item(118, 9).Value = 40000  # I118: was 0
$ws.Cells.Item(118, 11).Value = 40000  # K118: was 0
$ws.Cells.Item(118, 13).Value = -38343  # M118: was None
$ws.Cells.Item(132, 8).Value = 143620.58  # H132: was 145668.03
$ws.Cells.Item(132, 9).Value = 174697.11  # I132: was 177756.72
$ws.Cells.Item(132, 11).Value = 524091.33  # K132: was 533270.16
$ws.Cells.Item(132, 13).Value = -521561.33  # M132: was -530740.16
$ws.Cells.Item(136, 8).Value = 3627  # H136: was 3702.2307
$ws.Cells.Item(136, 9).Value = 2261.6365  # I136: was 2445.6191
$ws.Cells.Item(136, 10).Value = 8633.333000000001  # J136: was 8980
$ws.Cells.Item(136, 11).Value = 6784.9095  # K136: was 7336.8573
$ws.Cells.Item(136, 12).Value = 25899.999  # L136: was 26940
$ws.Cells.Item(136, 13).Value = -4234.9095  # M136: was -4786.8573
$ws.Cells.Item(136, 14).Value = -30999.999  # N136: was -32040

Write-Host "Applied 211 cell updates"